$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 8207.786
$ws.Range("I6").Value = 909.1667
$ws.Range("J6").Value = 51999.5
$ws.Range("K6").Value = 2727.5001
$ws.Range("L6").Value = 155998.5
$ws.Range("M6").Value = -2615.5001
$ws.Range("N6").Value = -156222.5

$ws.Range("H8").Value = 57.81818
$ws.Range("I8").Value = 57.81818
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 173.45454
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -34.45454000000001
$ws.Range("N8").ClearContents()

$ws.Range("H64").Value = 33870.406
$ws.Range("I64").Value = 69016.664
$ws.Range("J64").Value = 2859
$ws.Range("K64").Value = 69016.664
$ws.Range("L64").Value = 2859
$ws.Range("M64").Value = -68768.664
$ws.Range("N64").Value = -3355

$ws.Range("H67").Value = 33870.406
$ws.Range("I67").Value = 69016.664
$ws.Range("J67").Value = 2859
$ws.Range("K67").Value = 69016.664
$ws.Range("L67").Value = 2859
$ws.Range("M67").Value = -68158.664
$ws.Range("N67").Value = -4575

$ws.Range("H94").Value = 4579.9165
$ws.Range("I94").Value = 4579.9165
$ws.Range("K94").Value = 4579.9165
$ws.Range("M94").Value = -4128.9165

$ws.Range("H96").Value = 40044700
$ws.Range("I96").Value = 2310.6875
$ws.Range("J96").Value = 111231170
$ws.Range("K96").Value = 6932.0625
$ws.Range("L96").Value = 333693510
$ws.Range("M96").Value = -5559.0625
$ws.Range("N96").Value = -333696256

$ws.Range("H99").Value = 1066
$ws.Range("I99").Value = 1220.3077
$ws.Range("J99").Value = 397.33334
$ws.Range("K99").Value = 3660.9231
$ws.Range("L99").Value = 1192.00002
$ws.Range("M99").Value = -2162.9231
$ws.Range("N99").Value = -4188.000019999999

$ws.Range("H100").Value = 1294.0714
$ws.Range("I100").Value = 1254.4546
$ws.Range("J100").Value = 1439.3334
$ws.Range("K100").Value = 1254.4546
$ws.Range("L100").Value = 1439.3334
$ws.Range("M100").Value = -713.4546
$ws.Range("N100").Value = -2521.3334

$ws.Range("H101").Value = 688.2381
$ws.Range("I101").Value = 319.18182
$ws.Range("J101").Value = 1094.2
$ws.Range("K101").Value = 957.54546
$ws.Range("L101").Value = 3282.6
$ws.Range("M101").Value = 664.45454
$ws.Range("N101").Value = -6526.6

$ws.Range("H113").Value = 2341.9412
$ws.Range("I113").Value = 2095.2942
$ws.Range("J113").Value = 2588.5881
$ws.Range("K113").Value = 2095.2942
$ws.Range("L113").Value = 2588.5881
$ws.Range("M113").Value = 1158.7058
$ws.Range("N113").Value = -9096.588100000001

$ws.Range("H117").Value = 34864
$ws.Range("J117").Value = 34864
$ws.Range("L117").Value = 34864
$ws.Range("N117").Value = -44042

$ws.Range("H132").Value = 42099.5
$ws.Range("I132").Value = 6791.5625
$ws.Range("J132").Value = 112715.375
$ws.Range("K132").Value = 20374.6875
$ws.Range("L132").Value = 338146.125
$ws.Range("M132").Value = -17844.6875
$ws.Range("N132").Value = -343206.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1896.0741
$ws.Range("I2").Value = 1952.1
$ws.Range("J2").Value = 1736
$ws.Range("K2").Value = 1952.1
$ws.Range("L2").Value = 1736
$ws.Range("M2").Value = -1839.1
$ws.Range("N2").Value = -1962

$ws.Range("H97").Value = 651.9032
$ws.Range("I97").Value = 359.25
$ws.Range("J97").Value = 1184
$ws.Range("K97").Value = 359.25
$ws.Range("L97").Value = 1184
$ws.Range("M97").Value = 136.75
$ws.Range("N97").Value = -2176

$ws.Range("H113").Value = 46386
$ws.Range("J113").Value = 46386
$ws.Range("L113").Value = 46386
$ws.Range("N113").Value = -55064

$ws.Range("H116").Value = 1896.0741
$ws.Range("I116").Value = 1952.1
$ws.Range("J116").Value = 1736
$ws.Range("K116").Value = 1952.1
$ws.Range("L116").Value = 1736
$ws.Range("M116").Value = 341.9000000000001
$ws.Range("N116").Value = -6324

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1896.0741
$ws.Range("I3").Value = 1952.1
$ws.Range("J3").Value = 1736
$ws.Range("K3").Value = 1952.1
$ws.Range("L3").Value = 1736
$ws.Range("M3").Value = -1838.1
$ws.Range("N3").Value = -1964

$ws.Range("H82").Value = 15657.125
$ws.Range("I82").Value = 2628.5
$ws.Range("J82").Value = 20000
$ws.Range("K82").Value = 2628.5
$ws.Range("L82").Value = 20000
$ws.Range("M82").Value = -2245.5
$ws.Range("N82").Value = -20766

$ws.Range("H85").Value = 15657.125
$ws.Range("I85").Value = 2628.5
$ws.Range("J85").Value = 20000
$ws.Range("K85").Value = 2628.5
$ws.Range("L85").Value = 20000
$ws.Range("M85").Value = -1302.5
$ws.Range("N85").Value = -22652

$ws.Range("H94").Value = 1451
$ws.Range("I94").Value = 1201.25
$ws.Range("J94").Value = 2450
$ws.Range("K94").Value = 1201.25
$ws.Range("L94").Value = 2450
$ws.Range("M94").Value = -750.25
$ws.Range("N94").Value = -3352

$ws.Range("H99").Value = 2244.7896
$ws.Range("I99").Value = 2120.6667
$ws.Range("K99").Value = 2120.6667
$ws.Range("M99").Value = -622.6667000000002

$ws.Range("H111").Value = 47694
$ws.Range("J111").Value = 47694
$ws.Range("L111").Value = 47694
$ws.Range("N111").Value = -55874

$ws.Range("H116").Value = 43499
$ws.Range("J116").Value = 43499
$ws.Range("L116").Value = 43499
$ws.Range("N116").Value = -52677

$ws.Range("H117").Value = 44998
$ws.Range("J117").Value = 44998
$ws.Range("L117").Value = 44998
$ws.Range("N117").Value = -54176

$ws.Range("H134").Value = 3743.5881
$ws.Range("I134").Value = 2530.7827
$ws.Range("K134").Value = 7592.348100000001
$ws.Range("M134").Value = -5057.348100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 70003
$ws.Range("J3").Value = 70003
$ws.Range("L3").Value = 70003
$ws.Range("N3").Value = -70229

$ws.Range("H52").Value = 57499.5
$ws.Range("J52").Value = 57499.5
$ws.Range("L52").Value = 57499.5
$ws.Range("N52").Value = -58087.5

$ws.Range("H58").Value = 1842.2903
$ws.Range("I58").Value = 1525.9259
$ws.Range("J58").Value = 3977.75
$ws.Range("K58").Value = 1525.9259
$ws.Range("L58").Value = 3977.75
$ws.Range("M58").Value = -1322.9259
$ws.Range("N58").Value = -4383.75

$ws.Range("H99").Value = 1694.8846
$ws.Range("I99").Value = 1518.7
$ws.Range("J99").Value = 1805
$ws.Range("K99").Value = 1518.7
$ws.Range("L99").Value = 1805
$ws.Range("M99").Value = -20.70000000000005
$ws.Range("N99").Value = -4801

$ws.Range("H104").Value = 29997.889
$ws.Range("J104").Value = 29997.889
$ws.Range("L104").Value = 29997.889
$ws.Range("N104").Value = -35239.889

$ws.Range("H109").Value = 27128.273
$ws.Range("J109").Value = 27128.273
$ws.Range("L109").Value = 27128.273
$ws.Range("N109").Value = -29208.273

$ws.Range("H116").Value = 42364.25
$ws.Range("J116").Value = 42364.25
$ws.Range("L116").Value = 42364.25
$ws.Range("N116").Value = -51542.25

$ws.Range("H119").Value = 41250.668
$ws.Range("J119").Value = 41250.668
$ws.Range("L119").Value = 41250.668
$ws.Range("N119").Value = -50926.668

$ws.Range("H120").Value = 32613.727
$ws.Range("J120").Value = 32613.727
$ws.Range("L120").Value = 32613.727
$ws.Range("N120").Value = -39871.727

$ws.Range("H126").Value = 1694.8846
$ws.Range("I126").Value = 1518.7
$ws.Range("J126").Value = 1805
$ws.Range("K126").Value = 4556.1
$ws.Range("L126").Value = 5415
$ws.Range("M126").Value = -2086.1
$ws.Range("N126").Value = -10355

$ws.Range("H136").Value = 1842.2903
$ws.Range("I136").Value = 1525.9259
$ws.Range("J136").Value = 3977.75
$ws.Range("K136").Value = 4577.7777
$ws.Range("L136").Value = 11933.25
$ws.Range("M136").Value = -2027.7777
$ws.Range("N136").Value = -17033.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 3361.4285
$ws.Range("I123").Value = 3910
$ws.Range("J123").Value = 2950
$ws.Range("K123").Value = 11730
$ws.Range("L123").Value = 8850
$ws.Range("M123").Value = -9280
$ws.Range("N123").Value = -13750

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1640
$ws.Range("I14").Value = 280
$ws.Range("J14").Value = 3000
$ws.Range("K14").Value = 280
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = -112
$ws.Range("N14").Value = -3336

$ws.Range("H97").Value = 3886.5557
$ws.Range("I97").Value = 3167.5
$ws.Range("J97").Value = 5324.6665
$ws.Range("K97").Value = 3167.5
$ws.Range("L97").Value = 5324.6665
$ws.Range("M97").Value = -2671.5
$ws.Range("N97").Value = -6316.6665

$ws.Range("H102").Value = 1662.6
$ws.Range("I102").Value = 1853.25
$ws.Range("J102").Value = 900
$ws.Range("K102").Value = 1853.25
$ws.Range("L102").Value = 900
$ws.Range("M102").Value = -231.25
$ws.Range("N102").Value = -4144

$ws.Range("H116").Value = 36665.332
$ws.Range("J116").Value = 36665.332
$ws.Range("L116").Value = 36665.332
$ws.Range("N116").Value = -45843.332

$ws.Range("H122").Value = 1032.25
$ws.Range("I122").Value = 871
$ws.Range("J122").Value = 1258
$ws.Range("K122").Value = 2613
$ws.Range("L122").Value = 3774
$ws.Range("M122").Value = -163
$ws.Range("N122").Value = -8674

$ws.Range("H132").Value = 4616.577
$ws.Range("I132").Value = 1601.5333
$ws.Range("J132").Value = 8728
$ws.Range("K132").Value = 4804.5999
$ws.Range("L132").Value = 26184
$ws.Range("M132").Value = -2274.5999
$ws.Range("N132").Value = -31244

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1605.7
$ws.Range("I93").Value = 1557.5714
$ws.Range("J93").Value = 1631.6154
$ws.Range("K93").Value = 1557.5714
$ws.Range("L93").Value = 1631.6154
$ws.Range("M93").Value = -309.5714
$ws.Range("N93").Value = -4127.6154

$ws.Range("H112").Value = 36691.332
$ws.Range("J112").Value = 36691.332
$ws.Range("L112").Value = 36691.332
$ws.Range("N112").Value = -39645.332

$ws.Range("H116").Value = 45668
$ws.Range("J116").Value = 45668
$ws.Range("L116").Value = 45668
$ws.Range("N116").Value = -54846

$ws.Range("H119").Value = 36206
$ws.Range("J119").Value = 36206
$ws.Range("L119").Value = 36206
$ws.Range("N119").Value = -45882

$ws.Range("H120").Value = 46390
$ws.Range("J120").Value = 46390
$ws.Range("L120").Value = 46390
$ws.Range("N120").Value = -56066

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 405.8125
$ws.Range("I100").Value = 405.8125
$ws.Range("K100").Value = 811.625
$ws.Range("M100").Value = -270.625
